$d = $word.ActiveDocument

# 1. "Your netid:" paragraph -> merge runs, add netid value, drop one tab stop
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Your netid:" + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + [char]9, $true, $false, $false, $false, $false, $true, 1, $false, "Your netid: ldn190002" + [char]9 + [char]9 + [char]9 + [char]9 + [char]9 + [char]9, 2)

# 2. "Project author(s) netid(s):" -- just removes spell-check markup, no visible text change

# 3. "tf-idf" paragraph -- just removes spell-check markup, no visible text change

$d.Save()
